# mk3-shield-PCBWay-BOM.xlsx update
# - "added next version and support material"
#   * Title text PCB size 60x52mm -> 60x60mm
#   * Row 11 (item 5, JP1) designator becomes "5**" and gets a matching
#     "**" footnote further down
#   * New row 12 (item 6, JP4/6 SamTec stacking header) added
#   * New "**" / "DO NOT ORDER, INCLUDE AND  INSTALL!!!" footnote row (14)
#   * New blank/support rows + a merged note block appended below the table
#   * View zoom bumped to 110% on all sheets, selection moved to G2

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Title text: PCB size 60x52mm -> 60x60mm -------------------------------
$ws1.Range("D2").Value = "MK3-SHIELD   BOM  (Bill of Materials)  PCB: 60x60mm"

# --- Row 11: existing JP1 line gets a "5**" designator ----------------------
$ws1.Range("A11").Value = "5**"
$ws1.Range("B11").Value = "JP1"
$ws1.Range("C11").Value = 0
$ws1.Range("D11").Value = "-"
$ws1.Range("E11").Value = "-"
$ws1.Range("F11").Value = "2x3 jumper"
$ws1.Range("G11").Value = "-"
$ws1.Range("H11").Value = "Through Hole"
$ws1.Range("I11").Value = "do not populate (development/testing only)"

# --- Row 12: new part - JP4/6 SamTec board-to-board stacking connector -----
$ws1.Range("A12").Value = "6**"
$ws1.Range("B12").Value = "JP4/6"
$ws1.Range("C12").Value = 0
$ws1.Range("D12").Value = "SamTec"
$ws1.Range("E12").Value = "ESQ-122-13-L-T"
$ws1.Range("F12").Value = "3row 44pos in two 2x 3x22 ESQ-122-13-L-T stacking board inter connect"
$ws1.Range("H12").Value = "Through Hole"
$ws1.Range("I12").Value = "do not populate (bottom side mount!!)"

# --- Row 14: legend / footnote for the "**" markers -------------------------
$ws1.Range("A14").Value = "**"
$ws1.Range("B14").Value = "DO NOT ORDER, INCLUDE AND  INSTALL!!!"

# --- Row 25 height nudges slightly (13.5 -> 13.8) once the sheet grows -----
$ws1.Rows.Item(25).RowHeight = 13.8

# --- New support / note block below the table (rows 27-35) -----------------
# Formats copied from existing cells that already carry the right style.
$ws1.Range("D2:F4").Copy()
$ws1.Range("D27:F29").PasteSpecial(-4122)
$ws1.Range("D27:F29").Merge()

$ws1.Range("C7").Copy()
$ws1.Range("C30").PasteSpecial(-4122)
$ws1.Range("C31").PasteSpecial(-4122)

$ws1.Range("B7:F7").Copy()
$ws1.Range("B33:F33").PasteSpecial(-4122)

$ws1.Rows.Item(27).RowHeight = 13.8
$ws1.Rows.Item(28).RowHeight = 13.8
$ws1.Rows.Item(29).RowHeight = 13.8
$ws1.Rows.Item(30).RowHeight = 13.8
$ws1.Rows.Item(31).RowHeight = 21
$ws1.Rows.Item(32).RowHeight = 13.8
$ws1.Rows.Item(33).RowHeight = 13.8
$ws1.Rows.Item(34).RowHeight = 13.8
$ws1.Rows.Item(35).RowHeight = 13.8

$excel.CutCopyMode = 0

# --- Zoom to 110% on every sheet (view only, affects all 3 sheets) ---------
foreach ($wsName in @("Sheet2", "Sheet3", "Sheet1")) {
    $wsZoom = $wb.Worksheets.Item($wsName)
    $wsZoom.Activate()
    $excel.ActiveWindow.Zoom = 110
}

# --- Selection moves from I10 to G2 on Sheet1 -------------------------------
$ws1.Activate()
$ws1.Range("G2").Select() | Out-Null
